$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20, pushing existing rows 20-24 down to 21-25.
# Use the same style as the rest of the table (row 19) so formatting is preserved.
$ws.Rows.Item(20).Insert()

# Copy the row formatting from the row above (row 19) onto the newly inserted row 20,
# so the date cell keeps the date number format used throughout the table.
# Limit the copy to the used columns (A:R) to avoid touching the entire 16384-column row.
$ws.Range("A19:R19").Copy()
$ws.Range("A20:R20").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the values for the new row 20
$ws.Cells.Item(20, 1).Value = 7
$ws.Cells.Item(20, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(20, 3).Value = "Ñuble"
$ws.Cells.Item(20, 4).Value = 44837
$ws.Cells.Item(20, 5).Value = 16
$ws.Cells.Item(20, 6).Value = 300000000
$ws.Cells.Item(20, 7).Value = "Espárragos"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 200
$ws.Cells.Item(20, 11).Value = 1800
$ws.Cells.Item(20, 12).Value = 2000
$ws.Cells.Item(20, 13).Value = 1900
$ws.Cells.Item(20, 14).Value = "`$/kilo"
$ws.Cells.Item(20, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(20, 16).Value = 1900
$ws.Cells.Item(20, 17).Value = 1
$ws.Cells.Item(20, 18).Value = "Hortaliza"
